$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "nomorhp"
$ws.Range("B1").Value = "nomorhpexpected"
$ws.Range("C1").Value = "varTahun"
$ws.Range("D1").Value = "alamat"
$ws.Range("E1").Value = "alamatexpected"
$ws.Range("F1").Value = "pekerjaan"
$ws.Range("G1").Value = "pekerjaanexpected"
$ws.Range("H1").Value = "expected "

# --- Row 2 ---
$ws.Range("A2").Value = 87830815038
$ws.Range("B2").Value = 87830815038
$ws.Range("C2").Value = 1997
$ws.Range("D2").Value = "jl nuri 30"
$ws.Range("E2").Value = "jl nuri 30"
$ws.Range("F2").Value = "dokter"
$ws.Range("G2").Value = "dokter"
$ws.Range("H2").Value = "passed"

# --- Row 3 (write F3/G3 "mahasiswa" before D3/E3 "jl gagak 50" so the
#     shared-string table indices end up in the same order as the
#     reference workbook) ---
$ws.Range("F3").Value = "mahasiswa"
$ws.Range("G3").Value = "mahasiswa"
$ws.Range("D3").Value = "jl gagak 50"
$ws.Range("E3").Value = "jl gagak 50"
$ws.Range("A3").Value = 81904067865
$ws.Range("B3").Value = 81904067865
$ws.Range("C3").Value = 1998
$ws.Range("H3").Value = "failed"

# --- Row 4 ---
$ws.Range("A4").Value = "087830815asd"
$ws.Range("B4").Value = 87830815
$ws.Range("C4").Value = 1998
$ws.Range("D4").Value = "jl merpati 12"
$ws.Range("E4").Value = "jl merpati 12"
$ws.Range("F4").Value = "guru"
$ws.Range("G4").Value = "guru"
$ws.Range("H4").Value = "failed"

# --- Column widths (best effort; this runtime quantizes ColumnWidth to
#     1/6-character steps, so we pick the closest reachable value) ---
$ws.Columns(3).ColumnWidth = 18.5
$ws.Columns(4).ColumnWidth = 17.333333333333332
$ws.Columns(7).ColumnWidth = 18.166666666666668
$ws.Columns(8).ColumnWidth = 17.666666666666668

# --- Selection ---
[void]$ws.Range("F5").Select()
